$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 19) for "Billion Dollar Loser"
$ws.Range("A19").Value = "Billion Dollar Loser"
$ws.Range("B19").Value = "Reeves Wiedeman"

# Copy the existing date cell style (from row 18) so the new date cells
# use the same number format as the rest of the column, then set the
# underlying serial date values directly (avoids adding a time-of-day
# fraction that Value assignment with a DateTime would introduce).
$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D18").Copy()
$ws.Range("D19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C19").Value = 44253
$ws.Range("D19").Value = 44255

$ws.Range("E19").Value = "business;startups;wework;venture capital;visionary"
$ws.Range("F19").Value = "Audio"
$ws.Range("G19").Value = "10 Hours 56 Mins"

# Update the active cell selection to reflect the next empty row
$ws.Range("A20").Select()
